$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "36.738.72"
$ws.Range("E2").Value = "  +1.70%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.024.97"
$ws.Range("E3").Value = "  +0.64%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "248.41"
$ws.Range("E5").Value = "  -1.38%  "

$ws.Range("E6").Value = "  -0.92%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "62.86"
$ws.Range("E7").Value = "  +0.98%  "

$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.391"
$ws.Range("E9").Value = "  +5.79%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "58.09"
$ws.Range("E10").Value = "  -1.26%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0794"
$ws.Range("E11").Value = "  +6.38%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.104"
$ws.Range("E12").Value = "  -1.07%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.898"
$ws.Range("E13").Value = "  -1.31%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.40"
$ws.Range("E14").Value = "  +20.13%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.35"
$ws.Range("E15").Value = "  -3.21%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.315.84"
$ws.Range("E16").Value = "  +0.41%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.54"
$ws.Range("E17").Value = "  +2.41%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.024.27"
$ws.Range("E18").Value = "  +0.63%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "36.691.89"
$ws.Range("E19").Value = "  +1.71%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "72.19"
$ws.Range("E20").Value = "  +0.32%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0883"
$ws.Range("E21").Value = "  +3.29%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.39"
$ws.Range("E22").Value = "  +2.48%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.62"
$ws.Range("E23").Value = "  +1.16%  "

$ws.Range("E24").Value = "  +0.12%  "

$ws.Range("E25").Value = "  -5.38%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.33"
$ws.Range("E26").Value = "  +1.23%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.87"
$ws.Range("E27").Value = "  +3.98%  "

$ws.Range("B28").Value = "Kaspa"
$ws.Range("C28").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.137"
$ws.Range("E28").Value = "  +26.09%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "20.41"
$ws.Range("E29").Value = "  +4.37%  "

$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "159.99"
$ws.Range("E30").Value = "  -2.52%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.121"
$ws.Range("E31").Value = "  +0.55%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.07"
$ws.Range("E32").Value = "  -0.46%  "

$ws.Range("E33").Value = "  -1.26%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0620"
$ws.Range("E34").Value = "  +2.97%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.49"
$ws.Range("E35").Value = "  -0.21%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.39"
$ws.Range("E36").Value = "  -4.10%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.39"
$ws.Range("E37").Value = "  +9.77%  "

$ws.Range("E39").Value = "  +1.31%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.28"
$ws.Range("E40").Value = "  +33.47%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0998"
$ws.Range("E41").Value = "  -3.71%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.25"
$ws.Range("E42").Value = "  +3.22%  "

$ws.Range("E43").Value = "  +2.07%  "

$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "17.11"
$ws.Range("E44").Value = "  +3.15%  "

$ws.Range("B45").Value = "ARBITRUM"
$ws.Range("C45").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.13"
$ws.Range("E45").Value = "  +0.99%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0215"
$ws.Range("E46").Value = "  -0.10%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "93.98"
$ws.Range("E47").Value = "  +0.41%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.67"
$ws.Range("E48").Value = "  -1.82%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.372.58"
$ws.Range("E49").Value = "  -3.49%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.90"
$ws.Range("E50").Value = "  -0.36%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.208.32"
$ws.Range("E51").Value = "  +0.58%  "
